$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'liegen',
    'arten',
    'sorgen',
    'töten',
    'zeugen',
    'sagen',
    'flüchten',
    'zünden',
    'danken',
    'helfen',
    'leiden',
    'trennen',
    'leugnen',
    'loben',
    'zeigen',
    'regeln',
    'tropfen',
    'knarren',
    'wehen',
    'sperren',
    'mögen',
    'sterben',
    'kranken',
    'melden',
    'reiten',
    'folgen',
    'liefern',
    'streichen',
    'ärgern',
    'duschen',
    'schulden',
    'weichen',
    'äußern',
    'binden',
    'bilden',
    'biegen',
    'bitten',
    'gründen',
    'wundern',
    'decken',
    'freuen',
    'stoßen',
    'heilen',
    'klagen',
    'erben',
    'drehen',
    'nutzen',
    'merken',
    'tollen',
    'sprengen',
    'schwören',
    'malen',
    'fühlen',
    'bluten',
    'schlucken',
    'fügen',
    'betteln',
    'scheitern',
    'wenden',
    'trauen',
    'regnen',
    'wollen',
    'schwingen',
    'spüren',
    'tauchen',
    'heulen',
    'zögern',
    'buchen',
    'enden',
    'gnaden',
    'knien',
    'feiern',
    'räumen',
    'heben',
    'kosten',
    'ehren',
    'pissen',
    'strahlen',
    'kichern',
    'kaufen',
    'bellen',
    'sichern',
    'bergen',
    'machen',
    'schätzen',
    'wachsen',
    'spinnen',
    'treiben',
    'planen',
    'teilen',
    'husten',
    'stecken',
    'lügen',
    'schenken',
    'lockern',
    'hauen',
    'ächzen',
    'fangen',
    'posten',
    'fahren',
    'leisten',
    'gelten',
    'lohnen',
    'stürzen',
    'brauchen',
    'klingen',
    'fließen',
    'fischen',
    'fallen',
    'fehlen',
    'backen',
    'fassen',
    'grüßen',
    'atmen',
    'schwächen',
    'ziehen',
    'rufen',
    'greifen',
    'kehren',
    'filmen',
    'achten',
    'nerven',
    'altern',
    'knurren',
    'lehnen',
    'öffnen',
    'bauen',
    'lesen',
    'stürmen',
    'dringen',
    'wirken',
    'quälen',
    'kümmern',
    'zielen',
    'schreiten',
    'bremsen',
    'ändern',
    'tragen',
    'scheinen',
    'saufen',
    'formen',
    'mauern',
    'runden',
    'klettern',
    'spenden',
    'geben',
    'siegen',
    'suchen',
    'werfen',
    'treten',
    'fesseln',
    'jagen',
    'platzen',
    'graben',
    'seufzen',
    'fordern',
    'schmecken',
    'warnen',
    'hören',
    'münzen',
    'irren',
    'werden',
    'lassen',
    'wüten',
    'stehlen',
    'reizen',
    'pflegen',
    'rasen',
    'jubeln',
    'stellen',
    'fällen',
    'pfeifen',
    'narren',
    'wehtun',
    'flehen',
    'pflanzen',
    'zeichnen',
    'eignen',
    'heißen',
    'sitzen',
    'boxen',
    'schrecken',
    'dienen',
    'spielen',
    'läuten',
    'scheiden',
    'sinken',
    'bieten',
    'führen',
    'stammen',
    'nähen',
    'plaudern'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}

